$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (including the date number format/style index) from G5 to G6
# before assigning values, so the new row's date cell reuses the existing style.
$ws.Range("G5").Copy()
$ws.Range("G6").PasteSpecial(-4122)

$ws.Cells.Item(6, 1).Value = 10051.719999999999
$ws.Cells.Item(6, 2).Value = 9928.61
$ws.Cells.Item(6, 3).Value = 19.36
$ws.Cells.Item(6, 4).Value = 19.12
$ws.Cells.Item(6, 5).Value = $true
$ws.Cells.Item(6, 6).Value = -1.24
$ws.Cells.Item(6, 7).Value = 42612.674791666665
$ws.Cells.Item(6, 8).Value = $true
